# Applies the "Grant merge Ch5 5.1" copy-edit pass described by the diff:
#  - modelling -> modeling (US spelling) in the handful of body paragraphs
#    that use it (NOT in the "Data Modelling" / "Dynamic Modelling" headings)
#  - "dynamic aspect of" -> "dynamic aspects of"
#  - Oxford commas added before "and" in several entity-attribute lists
#  - "a order detail" -> "an order detail"
#  - a batch of small wording/grammar tweaks across the dynamic-modelling
#    narrative paragraphs (state/activity/sequence diagram sections)
#
# Each replacement is scoped to the specific paragraph's Range so that
# look-alike text elsewhere in the document (e.g. the section headings)
# is left untouched.

$d = $word.ActiveDocument

function Replace-InParagraph {
    param(
        [int]$Index,
        [string]$OldText,
        [string]$NewText
    )
    $rng = $d.Paragraphs($Index).Range
    $ok = $rng.Find.Execute($OldText, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $NewText, 2)
    if (-not $ok) {
        Write-Output "WARNING: paragraph $Index replace failed: '$OldText' -> '$NewText'"
    }
}

# --- Paragraph 2: "In this chapter ..." intro paragraph -------------------
Replace-InParagraph 2 "data modelling and dynamic modelling" "data modeling and dynamic modeling"
Replace-InParagraph 2 "In the data modelling section" "In the data modeling section"
Replace-InParagraph 2 "In the dynamic modelling section" "In the dynamic modeling section"
Replace-InParagraph 2 "dynamic aspect of our system" "dynamic aspects of our system"

# --- Paragraph 9: Customer entity description ------------------------------
Replace-InParagraph 9 "password and shipping addresses." "password, and shipping addresses."

# --- Paragraph 10: Product entity description -------------------------------
Replace-InParagraph 10 "brand and photo." "brand, and photo."

# --- Paragraph 11: Shopping Cart entity description -------------------------
Replace-InParagraph 11 "unit price and date when it is added." "unit price, and date when it is added."

# --- Paragraph 12: Purchase Order entity description -------------------------
Replace-InParagraph 12 "order status and customer Id." "order status, and customer Id."

# --- Paragraph 13: Order Detail entity description ---------------------------
Replace-InParagraph 13 "specific product and the total price" "specific product, and the total price"

# --- Paragraph 14: Preference entity description (red text) ------------------
Replace-InParagraph 14 "product Id number and the identifier" "product Id number, and the identifier"

# --- Paragraph 34: "order detail" relationship sentence -----------------------
Replace-InParagraph 34 "But a order detail" "But an order detail"

# --- Paragraph 67: Dynamic modelling section intro ----------------------------
Replace-InParagraph 67 "The dynamic modelling section" "The dynamic modeling section"

# --- Paragraph 69: State diagram narrative ------------------------------------
Replace-InParagraph 69 "state changes as the interaction between customer and vendor" "state changes with the interaction between the customer and vendor"
Replace-InParagraph 69 "After that, vendor can change the state" "After that, the vendor can change the state"
Replace-InParagraph 69 "if there is stock later, vendor can change the state" "if there is stock later, the vendor can change the state"

# --- Paragraph 75: Activity diagram (not logged in) narrative -----------------
Replace-InParagraph 75 "(include the new users)" "(including the new users)"
Replace-InParagraph 75 "such as home page, search page" "such as the home page, search page"
Replace-InParagraph 75 "go to home page to browse product list" "go to the home page to browse the product list"
Replace-InParagraph 75 "search the products by keywords in the search page" "search the products by keywords on the search page"
Replace-InParagraph 75 "filter the products by brands in the search page." "filter the products by brands on the search page."

# --- Paragraph 79: Activity diagram (logged in) narrative ----------------------
Replace-InParagraph 79 "customers can login to the system" "customers can log in to the system"
Replace-InParagraph 79 "create their own accounts in the register page, after registration" "create their accounts on the register page, and after registration"
Replace-InParagraph 79 "add the products which they want to buy" "add the products they want to buy"
Replace-InParagraph 79 "In the confirm page, customers are allowed" "On the confirm page, customers are allowed"
Replace-InParagraph 79 "passwords in the account management page" "passwords on the account management page"
Replace-InParagraph 79 ", order list page and account page" ", order list page, and account page"

# --- Paragraph 83: Activity diagram (vendor) narrative --------------------------
Replace-InParagraph 83 "In this diagram, vendor can enter the detail information of the new products" "In this diagram, the vendor can enter detailed information about the new products"
Replace-InParagraph 83 "When browsing the product list, vendor can search" "When browsing the product list, the vendor can search"
Replace-InParagraph 83 "or search the specific order by entering the order ID." "or search for the specific order by entering the order ID."
Replace-InParagraph 83 "In the order detail page, vendor can change the order states." "On the order detail page, the vendor can change the order states."
$quote = [char]8220
$unquote = [char]8221
$oldPending = "if the orders are in " + $quote + "pending" + $unquote + " state."
$newPending = "if the orders are in the " + $quote + "pending" + $unquote + " state."
Replace-InParagraph 83 $oldPending $newPending

# --- Paragraph 87: Sequence diagram intro ---------------------------------------
Replace-InParagraph 87 "exchange over time during the interaction" "exchange overtime during the interaction"

# --- Paragraph 88: Customer sequence diagram narrative --------------------------
Replace-InParagraph 88 "app from login to checkout process" "app from the login to the checkout process"
Replace-InParagraph 88 "their login credentials and the new users" "their login credentials, and the new users"
Replace-InParagraph 88 "product detail page and user can choose" "product detail page and the user can choose"

Write-Output "done"
